$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 0.1472031403336604
$ws.Range("E2").Value = 0.06105006105006105
$ws.Range("F2").Value = 0.4931506849315068
$ws.Range("G2").Value = 0.1372683596431023
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.1055408970976253

$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 0.2196595277320154
$ws.Range("G3").Value = 4.928131416837783
$ws.Range("H3").Value = 0.1021450459652707
$ws.Range("I3").Value = 0.2637130801687764

$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 0.2195389681668496
$ws.Range("G4").Value = 0.06863417982155114
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100

$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 0.05120327700972863
$ws.Range("I5").Value = 0.0526592943654555

$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100

$ws.Range("B7").Value = 100
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = 0.09818360333824251
$ws.Range("E7").Value = 0.1830384380719951
$ws.Range("F7").Value = 0.4390779363336992
$ws.Range("G7").Value = 1.371742112482853
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100

$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 0.1964636542239686
$ws.Range("E8").Value = 0.1830384380719951
$ws.Range("F8").Value = 0.4390779363336992
$ws.Range("G8").Value = 1.510989010989011
$ws.Range("H8").Value = 0.05112474437627813
$ws.Range("I8").Value = 0.05260389268805891

$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 0.1472754050073638
$ws.Range("E9").Value = 0.1217285453438831
$ws.Range("F9").Value = 0.5491488193300385
$ws.Range("G9").Value = 0.823045267489712
$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 100

$ws.Range("B10").Value = 100
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 0.04906771344455348
$ws.Range("E10").Value = 0.06086427267194157
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = 100
$ws.Range("I10").Value = 100

$ws.Range("B11").Value = 100
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 0.04911591355599214
$ws.Range("E11").Value = 0.06090133982947624
$ws.Range("F11").Value = 0.3846153846153846
$ws.Range("G11").Value = 1.241379310344828
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100

